$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.432.24'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '2.894.42'
$ws.Range('E3').Value = '  -2.94%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.89'
$ws.Range('E5').Value = '  -3.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.36'
$ws.Range('E6').Value = '  -6.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.552'
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').Value = '2.898.44'
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('E10').Value = '  -3.59%  '
$ws.Range('E11').Value = '  -4.63%  '
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').Value = '3.394.07'
$ws.Range('E13').Value = '  -3.07%  '
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('D15').Value = '60.473.34'
$ws.Range('E15').Value = '  -2.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.54'
$ws.Range('E16').Value = '  -4.18%  '
$ws.Range('D17').Value = '2.888.73'
$ws.Range('E17').Value = '  -3.29%  '
$ws.Range('E18').Value = '  -4.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.94'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.57'
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '362.11'
$ws.Range('E21').Value = '  -7.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.57'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.28'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').Value = '3.003.58'
$ws.Range('E25').Value = '  -4.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.449'
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.82'
$ws.Range('E29').Value = '  -6.85%  '
$ws.Range('D30').Value = '0.0₃0855'
$ws.Range('E30').Value = '  -8.54%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.44'
$ws.Range('E33').Value = '  -4.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.60'
$ws.Range('E34').Value = '  -6.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.33'
$ws.Range('E35').Value = '  -6.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.56'
$ws.Range('E36').Value = '  -6.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.997'
$ws.Range('E37').Value = '  -6.47%  '
$ws.Range('E38').Value = '  -5.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.81'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('D41').Value = '2.330.63'
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('E42').Value = '  -5.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.642'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.66'
$ws.Range('E44').Value = '  -6.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0572'
$ws.Range('E45').Value = '  -3.15%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.08'
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.998'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0234'
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.33'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0929'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '249.59'
$ws.Range('E51').Value = '  -4.71%  '
